$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: GroupID 100 -> 101 (VolID, PrimaryVolID, Admin unchanged)
$ws.Range("B2").Value = 101

# Row 3: GroupID 108 -> 100, VolID 589178b4... -> 293fe520...
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = "293fe520-7e35-444a-8955-f02a911fed1c"

# Row 4: VolID 293fe520... -> 34fb4310..., Admin 0 -> 1
$ws.Range("C4").Value = "34fb4310-9790-4b80-84cc-8c899f0308f7"
$ws.Range("E4").Value = 1

# Row 5: VolID 34fb4310... -> ec311095...
$ws.Range("C5").Value = "ec311095-16c4-4ea1-a9bc-9ddcda3b9b62"

# Row 10: VolID 46e0eab8... -> 589178b4..., PrimaryVolID 1 -> 0, Admin 1 -> 0
$ws.Range("C10").Value = "589178b4-aa4c-4276-a516-9460fa7714d3"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0

# Row 11: GroupID 100 -> 108, VolID ec311095... -> 46e0eab8..., PrimaryVolID 0 -> 1
$ws.Range("B11").Value = 108
$ws.Range("C11").Value = "46e0eab8-9d77-4a4d-a642-bed325a80ba2"
$ws.Range("D11").Value = 1

# Update active selection to C6 (from C15)
$ws.Range("C6").Select()
